$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell "Save" in H1 - copy the formatting (style) from the
# neighboring header cell G1 so it matches the rest of the header row,
# then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Data values for the new "Save" column (H2:H12)
$saveValues = @(0, 1, 0, 0, 1, 0, 0, 0, 1, 1, 0)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
